$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# The sheet is protected; unprotect to make edits, then re-protect at the end.
$ws.Unprotect()

# Update version number
$ws.Range("A2").Value = "Version 1.2.3"

# Update the existing row-5 instruction text to the new "no blank rows" instruction
$ws.Range("A5").Value = "Please use consecutive rows (no blank rows)."

# Insert a new row at row 6 (shifts everything from row 6 down by one)
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new instruction
$ws.Range("A6").Value = "Do not edit the header row of the 'Antibodies' sheet."

# Re-protect the sheet
$ws.Protect()
